$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 276, shifting rows 276:367 down to 277:368
$ws.Rows("276:276").Insert()

# Populate the new row 276 with the new data
$ws.Cells.Item(276, 1).Value = 4
$ws.Cells.Item(276, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(276, 3).Value = "Los Lagos"
$ws.Cells.Item(276, 4).Value = 44588
$ws.Cells.Item(276, 5).Value = 10
$ws.Cells.Item(276, 6).Value = "Fruta"
$ws.Cells.Item(276, 7).Value = 100108
$ws.Cells.Item(276, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(276, 9).Value = 100108006
$ws.Cells.Item(276, 10).Value = "Plátano"
$ws.Cells.Item(276, 11).Value = "Sin especificar"
$ws.Cells.Item(276, 12).Value = "Primera Pintón"
$ws.Cells.Item(276, 13).Value = 800
$ws.Cells.Item(276, 14).Value = 18000
$ws.Cells.Item(276, 15).Value = 19000
$ws.Cells.Item(276, 16).Value = 18500
$ws.Cells.Item(276, 17).Value = "`$/caja 20 kilos"
$ws.Cells.Item(276, 18).Value = "Ecuador"
$ws.Cells.Item(276, 19).Value = 925
$ws.Cells.Item(276, 20).Value = 20
